# Update column F (dSF) values for several rows to reflect the
# repulled/recalculated data, per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F19").Value = -2
$ws.Range("F25").Value = 3
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = -1
$ws.Range("F53").Value = -1
$ws.Range("F59").Value = -1
$ws.Range("F61").Value = -4
$ws.Range("F62").Value = -5
$ws.Range("F72").Value = -8
$ws.Range("F73").Value = 1
$ws.Range("F78").Value = -1
$ws.Range("F85").Value = 0
$ws.Range("F89").Value = 1
$ws.Range("F91").Value = 0
